# Scheduled-runner refresh of the Cerberus profit sheets: recomputed
# currentAveragePrice(NQ/HQ) / LevePrice(NQ/HQ) / LeveProfit(NQ/HQ) columns
# (H:N) for the leve rows whose underlying market data moved, across all
# eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 460.5
$ws.Cells.Item(2, 9).Value = 242.66667
$ws.Cells.Item(2, 10).Value = 678.3333
$ws.Cells.Item(2, 11).Value = 242.66667
$ws.Cells.Item(2, 12).Value = 678.3333
$ws.Cells.Item(2, 13).Value = -129.66667
$ws.Cells.Item(2, 14).Value = -904.3333
$ws.Cells.Item(48, 8).Value = 4579
$ws.Cells.Item(48, 9).Value = 4615
$ws.Cells.Item(48, 10).Value = 4573.857
$ws.Cells.Item(48, 11).Value = 13845
$ws.Cells.Item(48, 12).Value = 13721.571
$ws.Cells.Item(48, 13).Value = -13553
$ws.Cells.Item(48, 14).Value = -14305.571
$ws.Cells.Item(56, 8).Value = 4579
$ws.Cells.Item(56, 9).Value = 4615
$ws.Cells.Item(56, 10).Value = 4573.857
$ws.Cells.Item(56, 11).Value = 13845
$ws.Cells.Item(56, 12).Value = 13721.571
$ws.Cells.Item(56, 13).Value = -13311
$ws.Cells.Item(56, 14).Value = -14789.571
$ws.Cells.Item(64, 8).Value = 7361.25
$ws.Cells.Item(64, 10).Value = 8000
$ws.Cells.Item(64, 12).Value = 8000
$ws.Cells.Item(64, 14).Value = -8496
$ws.Cells.Item(67, 8).Value = 7361.25
$ws.Cells.Item(67, 10).Value = 8000
$ws.Cells.Item(67, 12).Value = 8000
$ws.Cells.Item(67, 14).Value = -9716
$ws.Cells.Item(74, 8).Value = 5963.364
$ws.Cells.Item(74, 9).Value = 5467.643
$ws.Cells.Item(74, 11).Value = 5467.643
$ws.Cells.Item(74, 13).Value = -4531.643
$ws.Cells.Item(77, 8).Value = 5963.364
$ws.Cells.Item(77, 9).Value = 5467.643
$ws.Cells.Item(77, 11).Value = 27338.215
$ws.Cells.Item(77, 13).Value = -22658.215
$ws.Cells.Item(98, 8).Value = 5884.0356
$ws.Cells.Item(98, 9).Value = 913.4211
$ws.Cells.Item(98, 10).Value = 16377.556
$ws.Cells.Item(98, 11).Value = 913.4211
$ws.Cells.Item(98, 12).Value = 16377.556
$ws.Cells.Item(98, 13).Value = 584.5789
$ws.Cells.Item(98, 14).Value = -19373.556
$ws.Cells.Item(122, 8).Value = 5884.0356
$ws.Cells.Item(122, 9).Value = 913.4211
$ws.Cells.Item(122, 10).Value = 16377.556
$ws.Cells.Item(122, 11).Value = 2740.2633
$ws.Cells.Item(122, 12).Value = 49132.66800000001
$ws.Cells.Item(122, 13).Value = -290.2633000000001
$ws.Cells.Item(122, 14).Value = -54032.66800000001
$ws.Cells.Item(131, 8).Value = 5062.5454
$ws.Cells.Item(131, 9).Value = 2876.4443
$ws.Cells.Item(131, 11).Value = 8629.332900000001
$ws.Cells.Item(131, 13).Value = -3589.332900000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2786.6865
$ws.Cells.Item(32, 9).Value = 2257.246
$ws.Cells.Item(32, 11).Value = 2257.246
$ws.Cells.Item(32, 13).Value = -1970.246
$ws.Cells.Item(61, 8).Value = 2450
$ws.Cells.Item(61, 9).Value = 1933.3334
$ws.Cells.Item(61, 10).Value = 4000
$ws.Cells.Item(61, 11).Value = 1933.3334
$ws.Cells.Item(61, 12).Value = 4000
$ws.Cells.Item(61, 13).Value = -1721.3334
$ws.Cells.Item(61, 14).Value = -4424
$ws.Cells.Item(136, 8).Value = 2450
$ws.Cells.Item(136, 9).Value = 1933.3334
$ws.Cells.Item(136, 10).Value = 4000
$ws.Cells.Item(136, 11).Value = 5800.0002
$ws.Cells.Item(136, 12).Value = 12000
$ws.Cells.Item(136, 13).Value = -3250.0002
$ws.Cells.Item(136, 14).Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 10840.158
$ws.Cells.Item(94, 9).Value = 353.9375
$ws.Cells.Item(94, 11).Value = 353.9375
$ws.Cells.Item(94, 13).Value = 97.0625
$ws.Cells.Item(105, 8).Value = 2728.926
$ws.Cells.Item(105, 9).Value = 1996.762
$ws.Cells.Item(105, 11).Value = 1996.762
$ws.Cells.Item(105, 13).Value = -249.7619999999999
$ws.Cells.Item(117, 8).Value = 189972
$ws.Cells.Item(117, 10).Value = 189972
$ws.Cells.Item(117, 12).Value = 189972
$ws.Cells.Item(117, 14).Value = -199150

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 5915.5
$ws.Cells.Item(62, 9).Value = 3332
$ws.Cells.Item(62, 10).Value = 8499
$ws.Cells.Item(62, 11).Value = 3332
$ws.Cells.Item(62, 12).Value = 8499
$ws.Cells.Item(62, 13).Value = -2708
$ws.Cells.Item(62, 14).Value = -9747
$ws.Cells.Item(65, 8).Value = 5915.5
$ws.Cells.Item(65, 9).Value = 3332
$ws.Cells.Item(65, 10).Value = 8499
$ws.Cells.Item(65, 11).Value = 16660
$ws.Cells.Item(65, 12).Value = 42495
$ws.Cells.Item(65, 13).Value = -13540
$ws.Cells.Item(65, 14).Value = -48735
$ws.Cells.Item(105, 8).Value = 1225.9
$ws.Cells.Item(105, 10).Value = 1272.5
$ws.Cells.Item(105, 12).Value = 1272.5
$ws.Cells.Item(105, 14).Value = -4766.5
$ws.Cells.Item(107, 8).Value = 534.625
$ws.Cells.Item(107, 9).Value = 326.23077
$ws.Cells.Item(107, 10).Value = 1437.6666
$ws.Cells.Item(107, 11).Value = 326.23077
$ws.Cells.Item(107, 12).Value = 1437.6666
$ws.Cells.Item(107, 13).Value = 1593.76923
$ws.Cells.Item(107, 14).Value = -5277.6666
$ws.Cells.Item(122, 8).Value = 2874.5833
$ws.Cells.Item(122, 9).Value = 2747.25
$ws.Cells.Item(122, 11).Value = 8241.75
$ws.Cells.Item(122, 13).Value = -5791.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(20, 8).Value = 1999.6666
$ws.Cells.Item(20, 10).Value = 1999.6666
$ws.Cells.Item(20, 12).Value = 5998.9998
$ws.Cells.Item(20, 14).Value = -6452.9998
$ws.Cells.Item(38, 8).Value = 311.7619
$ws.Cells.Item(38, 9).Value = 260.45456
$ws.Cells.Item(38, 11).Value = 781.36368
$ws.Cells.Item(38, 13).Value = -434.36368
$ws.Cells.Item(131, 8).Value = 20259620
$ws.Cells.Item(131, 10).Value = 22306222
$ws.Cells.Item(131, 12).Value = 66918666
$ws.Cells.Item(131, 14).Value = -66928746
$ws.Cells.Item(132, 8).Value = 1879.5143
$ws.Cells.Item(132, 10).Value = 2850
$ws.Cells.Item(132, 12).Value = 25650
$ws.Cells.Item(132, 14).Value = -30710

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 32495
$ws.Cells.Item(40, 10).Value = 32495
$ws.Cells.Item(40, 12).Value = 32495
$ws.Cells.Item(40, 14).Value = -32797
$ws.Cells.Item(98, 8).Value = 31708.25
$ws.Cells.Item(98, 10).Value = 31708.25
$ws.Cells.Item(98, 12).Value = 31708.25
$ws.Cells.Item(98, 14).Value = -37698.25
$ws.Cells.Item(122, 8).Value = 3595.6177
$ws.Cells.Item(122, 9).Value = 2811.96
$ws.Cells.Item(122, 10).Value = 5772.4443
$ws.Cells.Item(122, 11).Value = 8435.880000000001
$ws.Cells.Item(122, 12).Value = 17317.3329
$ws.Cells.Item(122, 13).Value = -5985.880000000001
$ws.Cells.Item(122, 14).Value = -22217.3329
$ws.Cells.Item(132, 8).Value = 3415.8096
$ws.Cells.Item(132, 9).Value = 2327.3333
$ws.Cells.Item(132, 11).Value = 6981.999899999999
$ws.Cells.Item(132, 13).Value = -4451.999899999999
$ws.Cells.Item(141, 8).Value = 90310
$ws.Cells.Item(141, 10).Value = 88970
$ws.Cells.Item(141, 12).Value = 88970
$ws.Cells.Item(141, 14).Value = -99330

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1271.6666
$ws.Cells.Item(22, 9).Value = 994.2143
$ws.Cells.Item(22, 10).Value = 1476.1052
$ws.Cells.Item(22, 11).Value = 994.2143
$ws.Cells.Item(22, 12).Value = 1476.1052
$ws.Cells.Item(22, 13).Value = -699.2143
$ws.Cells.Item(22, 14).Value = -2066.1052
$ws.Cells.Item(27, 8).Value = 1271.6666
$ws.Cells.Item(27, 9).Value = 994.2143
$ws.Cells.Item(27, 10).Value = 1476.1052
$ws.Cells.Item(27, 11).Value = 994.2143
$ws.Cells.Item(27, 12).Value = 1476.1052
$ws.Cells.Item(27, 13).Value = -887.2143
$ws.Cells.Item(27, 14).Value = -1690.1052
$ws.Cells.Item(68, 8).Value = 1903.6666
$ws.Cells.Item(68, 9).Value = 1615.5
$ws.Cells.Item(68, 11).Value = 1615.5
$ws.Cells.Item(68, 13).Value = -866.5
$ws.Cells.Item(71, 8).Value = 1903.6666
$ws.Cells.Item(71, 9).Value = 1615.5
$ws.Cells.Item(71, 11).Value = 8077.5
$ws.Cells.Item(71, 13).Value = -4333.5
$ws.Cells.Item(93, 8).Value = 1925.875
$ws.Cells.Item(93, 9).Value = 1710.8182
$ws.Cells.Item(93, 11).Value = 1710.8182
$ws.Cells.Item(93, 13).Value = -462.8181999999999
$ws.Cells.Item(100, 8).Value = 2761.5454
$ws.Cells.Item(100, 9).Value = 2080
$ws.Cells.Item(100, 11).Value = 2080
$ws.Cells.Item(100, 13).Value = -1539
$ws.Cells.Item(104, 8).Value = 74050.86
$ws.Cells.Item(104, 10).Value = 74050.86
$ws.Cells.Item(104, 12).Value = 74050.86
$ws.Cells.Item(104, 14).Value = -81038.86
$ws.Cells.Item(132, 8).Value = 2918.625
$ws.Cells.Item(132, 9).Value = 2328.9678
$ws.Cells.Item(132, 10).Value = 3993.8823
$ws.Cells.Item(132, 11).Value = 6986.903399999999
$ws.Cells.Item(132, 12).Value = 11981.6469
$ws.Cells.Item(132, 13).Value = -4456.903399999999
$ws.Cells.Item(132, 14).Value = -17041.6469

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(141, 8).Value = 73917.06
$ws.Cells.Item(141, 10).Value = 73917.06
$ws.Cells.Item(141, 12).Value = 73917.06
$ws.Cells.Item(141, 14).Value = -84277.06
